$d = $word.ActiveDocument

# 1) Shorten the "Attaching packages" console banner line
$d.Content.Find.Execute(
    "## ── Attaching packages ────────────────────────────────────────────────────────────────────────────────────────── tidyverse 1.2.1 ──",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## ── Attaching packages ─────────────────────────────────────── tidyverse 1.2.1 ──",
    2
)

# 2) Shorten the "Conflicts" console banner line
$d.Content.Find.Execute(
    "## ── Conflicts ───────────────────────────────────────────────────────────────────────────────────────────── tidyverse_conflicts() ──",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "## ── Conflicts ────────────────────────────────────────── tidyverse_conflicts() ──",
    2
)

# 3) Extend the figure caption with an extra sentence about removing grid lines
$d.Content.Find.Execute(
    "This figure includes the r values. These values can be removed by elimiating the lab = TRUE assignment",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This figure includes the r values. These values can be removed by elimiating the lab = TRUE assignment. TO remove the grid lines change ggtheme to theme_classic",
    2
)
